$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated results for case with 380 kV (vm_pu.xlsx)
$updates = @{
    "B2" = 1.02
    "C2" = 1.084580244026112
    "D2" = 1.03312030727396
    "E2" = 1.085074576338601
    "F2" = 1.091372745410925
    "I2" = 1.035677275773746
    "J2" = 1.089439326073989
    "K2" = 1.035923571682
    "L2" = 1.08773594734803
    "M2" = 1.094017906505344
    "N2" = 1.090986454877451
    "B3" = 1.02
    "C3" = 1.088554220956089
    "D3" = 1.033712597582539
    "E3" = 1.088594076997211
    "F3" = 1.094907733995827
    "I3" = 1.035967443082695
    "J3" = 1.093062457141762
    "K3" = 1.036325355068198
    "L3" = 1.091068096554292
    "M3" = 1.097366751325449
    "N3" = 1.09461473120692
    "B4" = 1.02
    "C4" = 1.091099880722971
    "D4" = 1.034092087236805
    "E4" = 1.090847350439783
    "F4" = 1.097169961554572
    "I4" = 1.036149929696968
    "J4" = 1.095381736569225
    "K4" = 1.03658088644504
    "L4" = 1.093200091228617
    "M4" = 1.099508478553141
    "N4" = 1.096937304277195
    "B5" = 1.02
    "C5" = 1.092164108474489
    "D5" = 1.034250736585222
    "E5" = 1.09178904086655
    "F5" = 1.098115160300037
    "I5" = 1.036225399593805
    "J5" = 1.096350932256881
    "K5" = 1.036687259228753
    "L5" = 1.094090778394488
    "M5" = 1.100403000555735
    "N5" = 1.097907876333895
    "B6" = 1.02
    "C6" = 1.092342452800936
    "D6" = 1.0342773227838
    "E6" = 1.091946832381985
    "F6" = 1.098273525854189
    "I6" = 1.03623799857866
    "J6" = 1.096513327957199
    "K6" = 1.036705058318497
    "L6" = 1.094240005038176
    "M6" = 1.100552856018284
    "N6" = 1.098070502654735
    "B7" = 1.02
    "C7" = 1.091114124171386
    "D7" = 1.034094210590322
    "E7" = 1.090859955067615
    "F7" = 1.097182614067086
    "I7" = 1.03615094301219
    "J7" = 1.09539470967109
    "K7" = 1.036582311922842
    "L7" = 1.093212014425301
    "M7" = 1.099520454000371
    "N7" = 1.096950295802352
    "B8" = 1.02
    "C8" = 1.085928727111042
    "D8" = 1.033321261131592
    "E8" = 1.086269103147404
    "F8" = 1.092572727252439
    "I8" = 1.035776439430744
    "J8" = 1.090669097336284
    "K8" = 1.036060285252374
    "L8" = 1.088867164676471
    "M8" = 1.095154987626646
    "N8" = 1.092217972555996
    "B9" = 1.02
    "C9" = 1.076584684746066
    "D9" = 1.031929828578845
    "E9" = 1.077986725810156
    "F9" = 1.084248648760612
    "I9" = 1.035075472029812
    "J9" = 1.08214085784312
    "K9" = 1.035105741197273
    "L9" = 1.081018157687755
    "M9" = 1.087261469720308
    "N9" = 1.083677621984548
    "B10" = 1.02
    "C10" = 1.070202816701081
    "D10" = 1.030981582126701
    "E10" = 1.072323574511657
    "F10" = 1.07855223380891
    "I10" = 1.034579580325208
    "J10" = 1.076307635532859
    "K10" = 1.034445207262348
    "L10" = 1.075644222990974
    "M10" = 1.081852383646115
    "N10" = 1.077836115829526
    "B11" = 1.02
    "C11" = 1.067400073508675
    "D11" = 1.030565885211432
    "E11" = 1.069834994285905
    "F11" = 1.076047959097567
    "I11" = 1.034357843132061
    "J11" = 1.073743830062632
    "K11" = 1.034153251532801
    "L11" = 1.073281037098021
    "M11" = 1.079472682495735
    "N11" = 1.07526866946146
    "B12" = 1.02
    "C12" = 1.066352816666969
    "D12" = 1.030410692213222
    "E12" = 1.068904907071711
    "F12" = 1.075111848652498
    "I12" = 1.03427440558879
    "J12" = 1.072785552156155
    "K12" = 1.034043895286085
    "L12" = 1.072397559058168
    "M12" = 1.078582875850084
    "N12" = 1.074309030690439
    "B13" = 1.02
    "C13" = 1.066577741719417
    "D13" = 1.030444017451906
    "E13" = 1.069104676798172
    "F13" = 1.075312919137159
    "I13" = 1.034292352193692
    "J13" = 1.072991380314381
    "K13" = 1.034067394098109
    "L13" = 1.072587329374588
    "M13" = 1.078774012451435
    "N13" = 1.074515151148259
    "B14" = 1.02
    "C14" = 1.067313634938957
    "D14" = 1.030553073026383
    "E14" = 1.069758231105423
    "F14" = 1.075970701982132
    "I14" = 1.034350968201893
    "J14" = 1.073664741768068
    "K14" = 1.034144230803
    "L14" = 1.073208125944328
    "M14" = 1.079399252174972
    "N14" = 1.075189468852446
    "B15" = 1.02
    "C15" = 1.06776621376657
    "D15" = 1.030620161223552
    "E15" = 1.070160142194451
    "F15" = 1.076375192770239
    "I15" = 1.034386940439408
    "J15" = 1.07407882340887
    "K15" = 1.034191451164332
    "L15" = 1.073589858408757
    "M15" = 1.079783696503827
    "N15" = 1.075604138536674
    "B16" = 1.02
    "C16" = 1.070387972720683
    "D16" = 1.031009061558281
    "E16" = 1.072487945480875
    "F16" = 1.078717619199092
    "I16" = 1.034594146780513
    "J16" = 1.076476964682975
    "K16" = 1.034464456695857
    "L16" = 1.075800276019674
    "M16" = 1.082009505496947
    "N16" = 1.078005685446458
    "B17" = 1.02
    "C17" = 1.072021802549288
    "D17" = 1.031251630327275
    "E17" = 1.073938197482961
    "F17" = 1.080176698508347
    "I17" = 1.034722229879739
    "J17" = 1.077970906458706
    "K17" = 1.034634102757089
    "L17" = 1.077176943806046
    "M17" = 1.083395481586696
    "N17" = 1.079501748790868
    "B18" = 1.02
    "C18" = 1.072971010524045
    "D18" = 1.03139262580981
    "E18" = 1.074780610409406
    "F18" = 1.081024135655641
    "I18" = 1.034796263411197
    "J18" = 1.078838650872932
    "K18" = 1.034732482586294
    "L18" = 1.077976451353359
    "M18" = 1.084200293586128
    "N18" = 1.080370725501694
    "B19" = 1.02
    "C19" = 1.073294033547345
    "D19" = 1.031440618926344
    "E19" = 1.075067266107048
    "F19" = 1.081312483337479
    "I19" = 1.034821393028898
    "J19" = 1.079133918488378
    "K19" = 1.03476593111551
    "L19" = 1.078248479744929
    "M19" = 1.084474109458128
    "N19" = 1.080666412431022
    "B20" = 1.02
    "C20" = 1.071846900769905
    "D20" = 1.031225655887883
    "E20" = 1.073782962391629
    "F20" = 1.080020529215615
    "I20" = 1.034708557759604
    "J20" = 1.077810999668594
    "K20" = 1.034615960619753
    "L20" = 1.077029601992348
    "M20" = 1.08324715399891
    "N20" = 1.079341614914773
    "B21" = 1.02
    "C21" = 1.067097106009377
    "D21" = 1.030520980690955
    "E21" = 1.069565935544834
    "F21" = 1.075777166563848
    "I21" = 1.034333737081951
    "J21" = 1.073466620450771
    "K21" = 1.034121629606209
    "L21" = 1.07302547576916
    "M21" = 1.079215298956916
    "N21" = 1.074991066180153
    "B22" = 1.02
    "C22" = 1.064074741281016
    "D22" = 1.030073373844432
    "E22" = 1.066881311074309
    "F22" = 1.073074863808862
    "I22" = 1.034091844788962
    "J22" = 1.070700479949765
    "K22" = 1.033805544157839
    "L22" = 1.070474900426695
    "M22" = 1.076646166572024
    "N22" = 1.072220997442355
    "B23" = 1.02
    "C23" = 1.065680461567747
    "D23" = 1.030311096351805
    "E23" = 1.068307715526734
    "F23" = 1.07451074588247
    "I23" = 1.034220674011907
    "J23" = 1.072170238255548
    "K23" = 1.033973613905696
    "L23" = 1.071830222177512
    "M23" = 1.078011432066209
    "N23" = 1.073692842973516
    "B24" = 1.02
    "C24" = 1.071925943009241
    "D24" = 1.031237394134564
    "E24" = 1.073853117229439
    "F24" = 1.080091106563909
    "I24" = 1.034714737687242
    "J24" = 1.077883265915003
    "K24" = 1.034624160036699
    "L24" = 1.077096190151747
    "M24" = 1.083314187973278
    "N24" = 1.079413983787542
    "B25" = 1.02
    "C25" = 1.079026199013511
    "D25" = 1.032293116840761
    "E25" = 1.080151961695951
    "F25" = 1.086425629927133
    "I25" = 1.035261651408421
    "J25" = 1.084370701988386
    "K25" = 1.035356708200856
    "L25" = 1.083071332386426
    "M25" = 1.08932711869216
    "N25" = 1.085910632764267
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
